# Auto-generated edit script: apply scheduled market-data refresh to Leve profit tables
# (columns H-N: currentAveragePrice*, LevePrice*, LeveProfit*) across all 8 job sheets.
$wb = $excel.ActiveWorkbook

# --- Sheet ALC: 38 cell updates ---
$ws = $wb.Worksheets.Item("ALC")
$updates = @(
    @("H40", "4350.769"),
    @("I40", "3378.875"),
    @("J40", "5905.8"),
    @("K40", "3378.875"),
    @("L40", "5905.8"),
    @("M40", "-3203.875"),
    @("N40", "-6255.8"),
    @("H43", "10078.625"),
    @("J43", "11435.615"),
    @("L43", "11435.615"),
    @("N43", "-11573.615"),
    @("H69", "4975.2"),
    @("I69", "4461.625"),
    @("K69", "13384.875"),
    @("M69", "-12510.875"),
    @("H72", "4975.2"),
    @("I72", "4461.625"),
    @("K72", "40154.625"),
    @("M72", "-35786.625"),
    @("H100", "9852.883"),
    @("J100", "15669.9"),
    @("L100", "15669.9"),
    @("N100", "-16751.9"),
    @("H132", "2160"),
    @("I132", "1224.8"),
    @("K132", "3674.4"),
    @("M132", "-1144.4"),
    @("H137", "3258.3684"),
    @("I137", "2679.2727"),
    @("K137", "8037.8181"),
    @("M137", "-5487.8181"),
    @("H141", "1490"),
    @("I141", "0"),
    @("J141", "1490"),
    @("K141", "0"),
    @("L141", "4470"),
    @("M141", $null),
    @("N141", "-14830"),
)
foreach ($u in $updates) {
    $ref = $u[0]
    $val = $u[1]
    if ($null -eq $val) {
        $ws.Range($ref).Value = ""
    } else {
        $ws.Range($ref).Value = [double]$val
    }
}

# --- Sheet ARM: 39 cell updates ---
$ws = $wb.Worksheets.Item("ARM")
$updates = @(
    @("H2", "15162796"),
    @("I2", "15162796"),
    @("K2", "15162796"),
    @("M2", "-15162683"),
    @("H49", "70000"),
    @("J49", "70000"),
    @("L49", "70000"),
    @("N49", "-70520"),
    @("H61", "12994.611"),
    @("I61", "12049.154"),
    @("K61", "12049.154"),
    @("M61", "-11837.154"),
    @("H74", "3504.0476"),
    @("I74", "1507.3334"),
    @("K74", "1507.3334"),
    @("M74", "-633.3334"),
    @("H77", "3504.0476"),
    @("I77", "1507.3334"),
    @("K77", "7536.666999999999"),
    @("M77", "-3168.666999999999"),
    @("H110", "10418422"),
    @("I110", "13889841"),
    @("J110", "4166.6665"),
    @("K110", "13889841"),
    @("L110", "4166.6665"),
    @("M110", "-13887796"),
    @("N110", "-8256.666499999999"),
    @("H116", "15162796"),
    @("I116", "15162796"),
    @("K116", "15162796"),
    @("M116", "-15160502"),
    @("H132", "4578.48"),
    @("I132", "3136"),
    @("K132", "9408"),
    @("M132", "-6878"),
    @("H136", "12994.611"),
    @("I136", "12049.154"),
    @("K136", "36147.462"),
    @("M136", "-33597.462"),
)
foreach ($u in $updates) {
    $ref = $u[0]
    $val = $u[1]
    if ($null -eq $val) {
        $ws.Range($ref).Value = ""
    } else {
        $ws.Range($ref).Value = [double]$val
    }
}

# --- Sheet BSM: 26 cell updates ---
$ws = $wb.Worksheets.Item("BSM")
$updates = @(
    @("H3", "15162796"),
    @("I3", "15162796"),
    @("K3", "15162796"),
    @("M3", "-15162682"),
    @("H74", "76633.336"),
    @("I74", "60000"),
    @("J74", "84950"),
    @("K74", "60000"),
    @("L74", "84950"),
    @("M74", "-59064"),
    @("N74", "-86822"),
    @("H77", "76633.336"),
    @("I77", "60000"),
    @("J77", "84950"),
    @("K77", "180000"),
    @("L77", "254850"),
    @("M77", "-175320"),
    @("N77", "-264210"),
    @("H81", "60779.8"),
    @("J81", "60779.8"),
    @("L81", "60779.8"),
    @("N81", "-62901.8"),
    @("H84", "60779.8"),
    @("J84", "60779.8"),
    @("L84", "182339.4"),
    @("N84", "-192947.4"),
)
foreach ($u in $updates) {
    $ref = $u[0]
    $val = $u[1]
    if ($null -eq $val) {
        $ws.Range($ref).Value = ""
    } else {
        $ws.Range($ref).Value = [double]$val
    }
}

# --- Sheet CRP: 56 cell updates ---
$ws = $wb.Worksheets.Item("CRP")
$updates = @(
    @("H16", "1392.3334"),
    @("I16", "1558.25"),
    @("J16", "1259.6"),
    @("K16", "1558.25"),
    @("L16", "1259.6"),
    @("M16", "-1271.25"),
    @("N16", "-1833.6"),
    @("H22", "718"),
    @("I22", "302"),
    @("K22", "302"),
    @("M22", "48"),
    @("H28", "20000"),
    @("J28", "20000"),
    @("L28", "20000"),
    @("N28", "-20490"),
    @("H58", "3499.8684"),
    @("I58", "2443.5"),
    @("J58", "4952.375"),
    @("K58", "2443.5"),
    @("L58", "4952.375"),
    @("M58", "-2240.5"),
    @("N58", "-5358.375"),
    @("H99", "4228"),
    @("I99", "2900"),
    @("K99", "2900"),
    @("M99", "-1402"),
    @("H107", "559.0476"),
    @("I107", "484.83334"),
    @("K107", "484.83334"),
    @("M107", "1435.16666"),
    @("H109", "39646.465"),
    @("J109", "39646.465"),
    @("L109", "39646.465"),
    @("N109", "-41726.465"),
    @("H113", "1392.3334"),
    @("I113", "1558.25"),
    @("J113", "1259.6"),
    @("K113", "1558.25"),
    @("L113", "1259.6"),
    @("M113", "611.75"),
    @("N113", "-5599.6"),
    @("H120", "39326"),
    @("I120", "0"),
    @("K120", "0"),
    @("M120", $null),
    @("H126", "4228"),
    @("I126", "2900"),
    @("K126", "8700"),
    @("M126", "-6230"),
    @("H136", "3499.8684"),
    @("I136", "2443.5"),
    @("J136", "4952.375"),
    @("K136", "7330.5"),
    @("L136", "14857.125"),
    @("M136", "-4780.5"),
    @("N136", "-19957.125"),
)
foreach ($u in $updates) {
    $ref = $u[0]
    $val = $u[1]
    if ($null -eq $val) {
        $ws.Range($ref).Value = ""
    } else {
        $ws.Range($ref).Value = [double]$val
    }
}

# --- Sheet CUL: 18 cell updates ---
$ws = $wb.Worksheets.Item("CUL")
$updates = @(
    @("H64", "2487.5"),
    @("I64", "1975"),
    @("J64", "3000"),
    @("K64", "5925"),
    @("L64", "9000"),
    @("M64", "-5655"),
    @("N64", "-9540"),
    @("H67", "2487.5"),
    @("I67", "1975"),
    @("J67", "3000"),
    @("K67", "5925"),
    @("L67", "9000"),
    @("M67", "-4989"),
    @("N67", "-10872"),
    @("H131", "8336095"),
    @("J131", "4873.263"),
    @("L131", "14619.789"),
    @("N131", "-24699.789"),
)
foreach ($u in $updates) {
    $ref = $u[0]
    $val = $u[1]
    if ($null -eq $val) {
        $ws.Range($ref).Value = ""
    } else {
        $ws.Range($ref).Value = [double]$val
    }
}

# --- Sheet GSM: 14 cell updates ---
$ws = $wb.Worksheets.Item("GSM")
$updates = @(
    @("H113", "61035.285"),
    @("I113", "83582.664"),
    @("J113", "44124.75"),
    @("K113", "83582.664"),
    @("L113", "44124.75"),
    @("M113", "-81412.664"),
    @("N113", "-48464.75"),
    @("H132", "5199.609"),
    @("I132", "3858.122"),
    @("J132", "16199.8"),
    @("K132", "11574.366"),
    @("L132", "48599.39999999999"),
    @("M132", "-9044.366"),
    @("N132", "-53659.39999999999"),
)
foreach ($u in $updates) {
    $ref = $u[0]
    $val = $u[1]
    if ($null -eq $val) {
        $ws.Range($ref).Value = ""
    } else {
        $ws.Range($ref).Value = [double]$val
    }
}

# --- Sheet LTW: 31 cell updates ---
$ws = $wb.Worksheets.Item("LTW")
$updates = @(
    @("H47", "30070"),
    @("J47", "0"),
    @("L47", "0"),
    @("N47", $null),
    @("H52", "30070"),
    @("J52", "0"),
    @("L52", "0"),
    @("N52", $null),
    @("H60", "29089"),
    @("I60", "29089"),
    @("J60", "0"),
    @("K60", "29089"),
    @("L60", "0"),
    @("M60", "-28580"),
    @("N60", $null),
    @("H74", "38650"),
    @("I74", "38650"),
    @("K74", "38650"),
    @("M74", "-37652"),
    @("H77", "38650"),
    @("I77", "38650"),
    @("K77", "115950"),
    @("M77", "-110958"),
    @("H98", "40355"),
    @("J98", "40355"),
    @("L98", "40355"),
    @("N98", "-46345"),
    @("H132", "3416.3225"),
    @("I132", "2734.875"),
    @("K132", "8204.625"),
    @("M132", "-5674.625"),
)
foreach ($u in $updates) {
    $ref = $u[0]
    $val = $u[1]
    if ($null -eq $val) {
        $ws.Range($ref).Value = ""
    } else {
        $ws.Range($ref).Value = [double]$val
    }
}

# --- Sheet WVR: 23 cell updates ---
$ws = $wb.Worksheets.Item("WVR")
$updates = @(
    @("H107", "2450.7144"),
    @("I107", "889.25"),
    @("J107", "4532.6665"),
    @("K107", "2667.75"),
    @("L107", "13597.9995"),
    @("M107", "-747.75"),
    @("N107", "-17437.9995"),
    @("H126", "6194.1333"),
    @("I126", "5911.4165"),
    @("K126", "17734.2495"),
    @("M126", "-15264.2495"),
    @("H131", "137521.42"),
    @("J131", "147000"),
    @("L131", "147000"),
    @("N131", "-157080"),
    @("H132", "4109.3423"),
    @("I132", "2751.6667"),
    @("K132", "8255.000100000001"),
    @("M132", "-5725.000100000001"),
    @("H136", "3275.5806"),
    @("I136", "1411.421"),
    @("K136", "4234.263"),
    @("M136", "-1684.263"),
)
foreach ($u in $updates) {
    $ref = $u[0]
    $val = $u[1]
    if ($null -eq $val) {
        $ws.Range($ref).Value = ""
    } else {
        $ws.Range($ref).Value = [double]$val
    }
}

Write-Output "Applied scheduled Sheets update across ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR."